$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '28.062.02'
$ws.Range('E2').Value = '  +0.00%  '

# Row 3
$ws.Range('D3').Value = '1.872.76'
$ws.Range('E3').Value = '  -0.88%  '

# Row 4
$ws.Range('D4').Value = '''1.005'
$ws.Range('E4').Value = '  +0.32%  '

# Row 5
$ws.Range('D5').Value = '''312.83'
$ws.Range('E5').Value = '  -0.37%  '

# Row 6
$ws.Range('E6').Value = '  +0.20%  '

# Row 7
$ws.Range('D7').Value = '''0.5145'
$ws.Range('E7').Value = '  +2.15%  '

# Row 8
$ws.Range('D8').Value = '''0.3820'

# Row 9
$ws.Range('D9').Value = '''0.08280'
$ws.Range('E9').Value = '  -10.30%  '

# Row 10
$ws.Range('D10').Value = '''1.113'
$ws.Range('E10').Value = '  -1.11%  '

# Row 11
$ws.Range('D11').Value = '''41.67'
$ws.Range('E11').Value = '  -0.31%  '

# Row 12
$ws.Range('D12').Value = '''6.222'
$ws.Range('E12').Value = '  -2.41%  '

# Row 13
$ws.Range('D13').Value = '1.877.37'
$ws.Range('E13').Value = '  -1.08%  '

# Row 14
$ws.Range('D14').Value = '''20.48'
$ws.Range('E14').Value = '  -1.45%  '

# Row 15
$ws.Range('D15').Value = '''7.200'

# Row 16
$ws.Range('D16').Value = '''1.005'
$ws.Range('E16').Value = '  +0.36%  '

# Row 17
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '''91.02'
$ws.Range('E17').Value = '  -1.27%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.00001094'
$ws.Range('E18').Value = '  -1.05%  '

# Row 19
$ws.Range('D19').Value = '''0.06651'
$ws.Range('E19').Value = '  -0.04%  '

# Row 20
$ws.Range('D20').Value = '''18.00'
$ws.Range('E20').Value = '  +1.01%  '

# Row 21
$ws.Range('E21').Value = '  +0.17%  '

# Row 22
$ws.Range('D22').Value = '''6.053'
$ws.Range('E22').Value = '  -2.43%  '

# Row 23
$ws.Range('D23').Value = '28.104.92'
$ws.Range('E23').Value = '  -0.04%  '

# Row 24
$ws.Range('D24').Value = '''11.16'
$ws.Range('E24').Value = '  -2.03%  '

# Row 25
$ws.Range('D25').Value = '''2.265'
$ws.Range('E25').Value = '  -2.29%  '

# Row 26
$ws.Range('D26').Value = '''2.582'
$ws.Range('E26').Value = '  +1.84%  '

# Row 27
$ws.Range('D27').Value = '2.097.15'
$ws.Range('E27').Value = '  -0.94%  '

# Row 28
$ws.Range('D28').Value = '''157.43'
$ws.Range('E28').Value = '  -0.61%  '

# Row 29
$ws.Range('D29').Value = '''20.59'

# Row 30
$ws.Range('D30').Value = '''125.84'
$ws.Range('E30').Value = '  -0.71%  '

# Row 31
$ws.Range('D31').Value = '''0.1060'
$ws.Range('E31').Value = '  +0.53%  '

# Row 32
$ws.Range('D32').Value = '''1.046'
$ws.Range('E32').Value = '  -2.67%  '

# Row 33
$ws.Range('D33').Value = '''5.608'
$ws.Range('E33').Value = '  +0.14%  '

# Row 34
$ws.Range('D34').Value = '''3.605'
$ws.Range('E34').Value = '  -0.46%  '

# Row 35
$ws.Range('D35').Value = '''9.641'
$ws.Range('E35').Value = '  +1.77%  '

# Row 36
$ws.Range('D36').Value = '''0.02455'
$ws.Range('E36').Value = '  +2.29%  '

# Row 37
$ws.Range('D37').Value = '''0.06567'
$ws.Range('E37').Value = '  -0.19%  '

# Row 38
$ws.Range('D38').Value = '''0.2166'
$ws.Range('E38').Value = '  -1.41%  '

# Row 39
$ws.Range('D39').Value = '''1.216'
$ws.Range('E39').Value = '  +0.23%  '

# Row 40
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '''1.247'
$ws.Range('E40').Value = '  -6.69%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.6450'
$ws.Range('E41').Value = '  +0.45%  '

# Row 42
$ws.Range('E42').Value = '  -2.07%  '

# Row 43
$ws.Range('D43').Value = '''4.879'
$ws.Range('E43').Value = '  -1.51%  '

# Row 44
$ws.Range('D44').Value = '''0.6099'
$ws.Range('E44').Value = '  +1.12%  '

# Row 45
$ws.Range('D45').Value = '''13.13'
$ws.Range('E45').Value = '  -1.00%  '

# Row 46
$ws.Range('D46').Value = '''1.298'
$ws.Range('E46').Value = '  -0.20%  '

# Row 47
$ws.Range('D47').Value = '''3.669'
$ws.Range('E47').Value = '  -0.41%  '

# Row 48
$ws.Range('D48').Value = '''2.009'
$ws.Range('E48').Value = '  +0.56%  '

# Row 49
$ws.Range('D49').Value = '''1.219'
$ws.Range('E49').Value = '  +2.32%  '

# Row 50
$ws.Range('D50').Value = '''120.82'
$ws.Range('E50').Value = '  -1.03%  '

# Row 51
$ws.Range('D51').Value = '''80.51'
$ws.Range('E51').Value = '  +2.01%  '
